$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-28 Tuesday" "2025-01-29 Wednesday"

Replace-Text "756×8=6048" "338×6=2028"
Replace-Text "864×4=3456" "973×3=2919"
Replace-Text "328×7=2296" "494×5=2470"
Replace-Text "692×4=2768" "920×6=5520"
Replace-Text "117×3=351" "283×3=849"
Replace-Text "413×6=2478" "571×5=2855"
Replace-Text "614×5=3070" "339×3=1017"
Replace-Text "349×5=1745" "673×9=6057"
Replace-Text "187×5=935" "135×3=405"
Replace-Text "371×8=2968" "869×5=4345"
Replace-Text "908×5=4540" "893×6=5358"
Replace-Text "482×8=3856" "383×9=3447"
Replace-Text "288×9=2592" "447×8=3576"
Replace-Text "405×9=3645" "805×5=4025"
Replace-Text "860×7=6020" "197×6=1182"
Replace-Text "218×9=1962" "490×6=2940"
Replace-Text "704×9=6336" "890×5=4450"
Replace-Text "908×6=5448" "293×3=879"
Replace-Text "910×9=8190" "836×5=4180"
Replace-Text "724×2=1448" "775×3=2325"
Replace-Text "147×6=882" "714×9=6426"
Replace-Text "415×8=3320" "657×9=5913"
Replace-Text "458×2=916" "904×5=4520"
Replace-Text "751×5=3755" "870×4=3480"
Replace-Text "370×4=1480" "140×6=840"

Write-Output "Done"
